# Add a new "2022-Q4" quarter sheet (copied from the "2022-Q3" layout) right
# after the "总计" (total) summary sheet, and record the new quarter in the
# summary table. All the other quarter sheets simply shift one tab position
# to the right - their own names/content are untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (same columns/
#    styles/fund rows), inserted right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $totalSheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Update the fund-position figures for the new quarter (these columns are
# stored as text in the workbook, so force a text number format first).
$q4Sheet.Range("D2:G3").NumberFormat = "@"

$q4Sheet.Range("D2").Value = "0.50"
$q4Sheet.Range("E2").Value = "90.73"
$q4Sheet.Range("F2").Value = "4.51"
$q4Sheet.Range("G2").Value = "0.0226"

$q4Sheet.Range("D3").Value = "0.09"
$q4Sheet.Range("E3").Value = "90.73"
$q4Sheet.Range("F3").Value = "4.51"
$q4Sheet.Range("G3").Value = "0.0041"

# ---------------------------------------------------------------------
# 2) Insert a new row at the top of the "总计" data table for 2022-Q4,
#    pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.03

# Clean up the row-2 formatting so it matches the other data rows: B:D plain,
# A styled like the index column in the rows below it.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
